$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force column D and E to Text format so that numeric-looking
# strings (e.g. "66.253.88", "584.50", "0.0000257") are preserved exactly
# as text instead of being auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '66.253.88'
$ws.Range('E2').Value = '  +6.70%  '
$ws.Range('D3').Value = '3.015.50'
$ws.Range('E3').Value = '  +3.65%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '584.50'
$ws.Range('E5').Value = '  +3.01%  '
$ws.Range('D6').Value = '163.14'
$ws.Range('E6').Value = '  +13.44%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').Value = '0.520'
$ws.Range('E8').Value = '  +3.85%  '
$ws.Range('D9').Value = '3.013.75'
$ws.Range('E9').Value = '  +3.67%  '
$ws.Range('D10').Value = '6.70'
$ws.Range('E10').Value = '  -3.99%  '
$ws.Range('E11').Value = '  +5.07%  '
$ws.Range('E12').Value = '  +5.91%  '
$ws.Range('D13').Value = '0.0000257'
$ws.Range('E13').Value = '  +7.99%  '
$ws.Range('D14').Value = '34.78'
$ws.Range('E14').Value = '  +6.96%  '
$ws.Range('E15').Value = '  -0.68%  '
$ws.Range('D16').Value = '66.272.48'
$ws.Range('E16').Value = '  +6.77%  '
$ws.Range('D17').Value = '3.521.01'
$ws.Range('E17').Value = '  +3.72%  '
$ws.Range('D18').Value = '6.95'
$ws.Range('E18').Value = '  +5.80%  '
$ws.Range('D19').Value = '3.018.43'
$ws.Range('E19').Value = '  +3.16%  '
$ws.Range('D20').Value = '457.39'
$ws.Range('E20').Value = '  +6.16%  '
$ws.Range('D21').Value = '14.00'
$ws.Range('E21').Value = '  +7.20%  '
$ws.Range('D22').Value = '0.689'
$ws.Range('E22').Value = '  +5.63%  '
$ws.Range('D23').Value = '7.39'
$ws.Range('E23').Value = '  +7.92%  '
$ws.Range('D24').Value = '82.49'
$ws.Range('E24').Value = '  +4.75%  '
$ws.Range('E25').Value = '  +15.44%  '
$ws.Range('E26').Value = '  +3.61%  '
$ws.Range('D27').Value = '10.56'
$ws.Range('E27').Value = '  +5.28%  '
$ws.Range('E28').Value = '  -0.07%  '
$ws.Range('D29').Value = '8.13'
$ws.Range('E29').Value = '  +16.73%  '
$ws.Range('E30').Value = '  +18.50%  '
$ws.Range('D31').Value = '0.0000106'
$ws.Range('E31').Value = '  -5.80%  '
$ws.Range('E32').Value = '  +4.44%  '
$ws.Range('D33').Value = '27.42'
$ws.Range('E33').Value = '  +6.87%  '
$ws.Range('D34').Value = '0.111'
$ws.Range('E34').Value = '  +5.55%  '
$ws.Range('E35').Value = '  -0.05%  '
$ws.Range('D36').Value = '0.994'
$ws.Range('E36').Value = '  +3.91%  '
$ws.Range('B37').Value = 'Stacks'
$ws.Range('C37').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D37').Value = '2.20'
$ws.Range('E37').Value = '  +16.33%  '
$ws.Range('B38').Value = 'Filecoin'
$ws.Range('C38').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D38').Value = '5.82'
$ws.Range('E38').Value = '  +8.12%  '
$ws.Range('D39').Value = '3.03'
$ws.Range('E39').Value = '  +3.97%  '
$ws.Range('D40').Value = '50.07'
$ws.Range('E40').Value = '  +2.43%  '
$ws.Range('E41').Value = '  +16.12%  '
$ws.Range('D42').Value = '0.124'
$ws.Range('E42').Value = '  +8.91%  '
$ws.Range('E43').Value = '  +7.04%  '
$ws.Range('E44').Value = '  +4.08%  '
$ws.Range('D45').Value = '400.74'
$ws.Range('E45').Value = '  +16.00%  '
$ws.Range('D46').Value = '0.0362'
$ws.Range('E46').Value = '  +7.31%  '
$ws.Range('D47').Value = '2.804.07'
$ws.Range('E47').Value = '  +2.93%  '
$ws.Range('D48').Value = '134.52'
$ws.Range('E48').Value = '  +0.93%  '
$ws.Range('E49').Value = '  +0.01%  '
$ws.Range('D50').Value = '23.97'
$ws.Range('E50').Value = '  +11.72%  '
$ws.Range('E51').Value = '  +4.83%  '

# Restore the default (Normal) style on the affected range so the cell
# styling matches the original workbook (no leftover explicit text format).
$ws.Range("D2:E51").Style = "Normal"
